# The commit (msg: "finish updating acc tests") extends the last bullet of
# the acceptance-test fixes list:
#   1. Gives that bullet real text ("...mise mivchanei kabala...").
#   2. Adds two blank spacer paragraphs after it.
#   3. Adds a final, un-styled paragraph containing a long
#      `-k "test_a or test_b or ..."` pytest filter expression, split across
#      many runs with <w:proofErr> spell-check markers (Word's own doing)
#      and a couple of runs in Arial. The pre-existing "_GoBack" bookmark
#      keeps its original position, which now falls in the middle of the
#      word "test_managerDoingThings" (between the "manage" run and the
#      "rDoingThings" run).
#
# Reproducing that exact run layout (and the bookmark sitting mid-word)
# through Find/Replace or TypeText would be unreliable, so instead we
# collapse a Range to the very end of the document and hand Word the
# replacement OOXML directly via Range.InsertXML -- the same mechanism
# Word itself uses under the hood when pasting/merging OOXML fragments.
# InsertXML on a Range collapsed at a paragraph's end swaps in the supplied
# <w:p>...</w:p> elements for that paragraph (additional sibling <w:p>
# elements become new paragraphs right after it), so the existing bookmark
# is simply included, verbatim, at the right spot inside the new XML.

$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00800D66" w:rsidRPr="00C958C5" w:rsidRDefault="00800D66" w:rsidP="00800D66">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:bidi/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>הוספת מבחני קבלה נוספים הבודקים מקרי כישלון רבים ואילוצי נכונות שהיו חסרים בגרסה הקודמת</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:bidi/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:bidi/>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial"/>
      <w:rtl/>
    </w:rPr>
    <w:t>-</w:t>
  </w:r>
  <w:r>
    <w:t>k "</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_setup</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_register</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_login</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_search</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_saveItem</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_watchCart</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_removeItemFromCart</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_changeItemQuantityInCart</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_logout</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_addStore</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_addItemToStore</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_removeItemFromStore</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_changeItemInStore</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_test_addOwner</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_removeOwner</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_addManager</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_removeManager</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>test_</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial"/>
    </w:rPr>
    <w:t>manage</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial"/>
    </w:rPr>
    <w:t>rDoingThings</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve"> or </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial"/>
    </w:rPr>
    <w:t>test_removeUser</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial"/>
      <w:rtl/>
    </w:rPr>
    <w:t>"</w:t>
  </w:r>
</w:p>
'@

$r.InsertXML($xml)
